$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update OCR "ground truth" TEXT column values (debug data) ---
$ws.Range("E4").Value = "重壤之光留顯號于千秋永誌外家之福欽"
$ws.Range("E5").Value = "聖祖皇帝五旬大慶節欽奉"
$ws.Range("E6").Value = "靈應肆丕承"

# --- Column C gets its own (slightly narrower) width, split off from C:D ---
$ws.Columns("C").ColumnWidth = 37.33

# --- View state: zoom to 96%, scroll so row 6 is at the top, select D17 ---
$excel.ActiveWindow.Zoom = 96
$excel.ActiveWindow.ScrollRow = 6
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("D17").Select()
